$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'329.03"
$ws.Range("E2").Value = "'0.63%"
$ws.Range("D3").Value = "'44.09"
$ws.Range("E3").Value = "'0.85%"
$ws.Range("D4").Value = "'5.577"
$ws.Range("E4").Value = "'1.62%"
$ws.Range("D5").Value = "'0.08066"
$ws.Range("E5").Value = "'0.02%"
$ws.Range("D6").Value = "'1.985"
$ws.Range("E6").Value = "'5.80%"
$ws.Range("B7").Value = "GateToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D7").Value = "'4.329"
$ws.Range("E7").Value = "'1.14%"
$ws.Range("B8").Value = "BTSEToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D8").Value = "'2.574"
$ws.Range("E8").Value = "'-5.44%"
$ws.Range("B9").Value = "MXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D9").Value = "'0.9530"
$ws.Range("E9").Value = "'1.80%"
$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D10").Value = "'0.1157"
$ws.Range("E10").Value = "'-0.66%"
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").Value = "'0.1857"
$ws.Range("E11").Value = "'-1.84%"
$ws.Range("B12").Value = "MCDex"
$ws.Range("C12").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D12").Value = "'11.87"
$ws.Range("E12").Value = "'39.14%"
$ws.Range("B13").Value = "MandalaExchangeToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D13").Value = "'0.09844"
$ws.Range("E13").Value = "'2.74%"
$ws.Range("B14").Value = "BitrueCoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D14").Value = "'0.04755"
$ws.Range("E14").Value = "'14.65%"
$ws.Range("B15").Value = "BitMartToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D15").Value = "'0.1069"
$ws.Range("E15").Value = "'0.36%"
$ws.Range("B16").Value = "BitForexToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D16").Value = "'0.001285"
$ws.Range("E16").Value = "'0.27%"
$ws.Range("B17").Value = "CoinExToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D17").Value = "'0.04241"
$ws.Range("E17").Value = "'-2.40%"
$ws.Range("B18").Value = "TigerCash"
$ws.Range("C18").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D18").Value = "'0.005888"
$ws.Range("E18").Value = "'-1.41%"
$ws.Range("B19").Value = "LEO"
$ws.Range("C19").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D19").Value = "'3.371"
$ws.Range("E19").Value = "'-5.44%"
$ws.Range("D21").Value = "'0.1410"
$ws.Range("E21").Value = "'3.34%"
$ws.Range("D22").Value = "'0.2509"
$ws.Range("E22").Value = "'-3.20%"
$ws.Range("D23").Value = "'0.001254"
$ws.Range("E23").Value = "'1.67%"
$ws.Range("D24").Value = "'0.004330"
$ws.Range("E24").Value = "'-0.37%"
$ws.Range("E25").Value = "'-3.43%"
$ws.Range("E26").Value = "'-0.44%"
$ws.Range("D38").Value = "'0.02635"
$ws.Range("E38").Value = "'-0.87%"
$ws.Range("D39").Value = "'0.05543"
$ws.Range("E39").Value = "'2.08%"
$ws.Range("E40").Value = "'-0.78%"
$ws.Range("D41").Value = "'0.1408"
$ws.Range("E41").Value = "'1.40%"
$ws.Range("D42").Value = "'0.008085"
$ws.Range("E42").Value = "'-29.42%"
$ws.Range("D43").Value = "'0.002018"
$ws.Range("E43").Value = "'-5.20%"
$ws.Range("D44").Value = "'0.008865"
$ws.Range("E44").Value = "'-8.42%"
$ws.Range("D45").Value = "'0.00007094"
$ws.Range("E45").Value = "'3.30%"
$ws.Range("E46").Value = "'-0.12%"
$ws.Range("E47").Value = "'1.22%"
$ws.Range("D48").Value = "'0.003576"
$ws.Range("E48").Value = "'0.24%"
$ws.Range("D49").Value = "'0.00002103"
$ws.Range("E49").Value = "'-0.12%"
$ws.Range("D50").Value = "'0.0002003"
$ws.Range("E50").Value = "'-0.12%"
